$d = $word.ActiveDocument

# Locate the paragraphs involved:
#  - $anchorPara: the paragraph ending "...other class materials." that should retain
#    the _GoBack bookmark at its very end.
#  - $startPara:  first paragraph of the block to remove ("Choose your group")
#  - $endPara:    last paragraph of the block to remove ("...its unstructured nature.")
$anchorPara = $null
$startPara = $null
$endPara = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*other class materials*") {
        $anchorPara = $p
    }
    if ($t -like "*Choose your group*") {
        $startPara = $p
    }
    if ($t -like "*its unstructured nature*") {
        $endPara = $p
    }
}

# Remove the whole block of paragraphs from "Choose your group" through the
# paragraph ending in "its unstructured nature." (inclusive).
$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()

# Re-fetch the anchor paragraph's current end position (positions shifted
# after the delete above).
$insertPos = $anchorPara.Range.End - 1

# Work around a boundary quirk: creating a zero-length Range exactly at
# (paragraph.End - 1) and adding a bookmark there can mis-place the
# bookmark. To avoid this, temporarily widen the paragraph by inserting a
# placeholder character after the target position, add the bookmark just
# before the placeholder (now a safe, non-boundary position), then remove
# the placeholder. Bookmarks remain anchored to their text position when
# surrounding text is edited.
$placeholderRange = $d.Range($insertPos, $insertPos)
$placeholderRange.InsertAfter("X")

$bmRange = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholderDeleteRange = $d.Range($insertPos, $insertPos + 1)
$placeholderDeleteRange.Delete()
